$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 6510
$ws.Range("J69").Value = 6510
$ws.Range("L69").Value = 19530
$ws.Range("N69").Value = -21278
$ws.Range("H70").Value = 1560
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H72").Value = 6510
$ws.Range("J72").Value = 6510
$ws.Range("L72").Value = 58590
$ws.Range("N72").Value = -67326
$ws.Range("H73").Value = 1560
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H98").Value = 2290.0588
$ws.Range("I98").Value = 1565.6428
$ws.Range("K98").Value = 1565.6428
$ws.Range("M98").Value = -67.64280000000008
$ws.Range("H116").Value = 2457.5557
$ws.Range("I116").Value = 2302
$ws.Range("J116").Value = 2652
$ws.Range("K116").Value = 2302
$ws.Range("L116").Value = 2652
$ws.Range("M116").Value = 1140
$ws.Range("N116").Value = -9536
$ws.Range("H122").Value = 2290.0588
$ws.Range("I122").Value = 1565.6428
$ws.Range("K122").Value = 4696.928400000001
$ws.Range("M122").Value = -2246.928400000001
$ws.Range("H129").Value = 768
$ws.Range("I129").Value = 466
$ws.Range("J129").Value = 1825
$ws.Range("K129").Value = 1398
$ws.Range("L129").Value = 5475
$ws.Range("M129").Value = 3602
$ws.Range("N129").Value = -15475
$ws.Range("H134").Value = 44046.188
$ws.Range("J134").Value = 44046.188
$ws.Range("L134").Value = 44046.188
$ws.Range("N134").Value = -54186.188
$ws.Range("H137").Value = 11906041
$ws.Range("I137").Value = 1406.9375
$ws.Range("J137").Value = 50000870
$ws.Range("K137").Value = 4220.8125
$ws.Range("L137").Value = 150002610
$ws.Range("M137").Value = -1670.8125
$ws.Range("N137").Value = -150007710
$ws.Range("H138").Value = 4062.13
$ws.Range("J138").Value = 4415.898
$ws.Range("L138").Value = 13247.694
$ws.Range("N138").Value = -23527.694

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 1669629
$ws.Range("J112").Value = 1669629
$ws.Range("L112").Value = 1669629
$ws.Range("N112").Value = -1672583

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2012.7858
$ws.Range("I86").Value = 1552.2307
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 1552.2307
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = -429.2307000000001
$ws.Range("N86").Value = -10246
$ws.Range("H89").Value = 2012.7858
$ws.Range("I89").Value = 1552.2307
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 7761.1535
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = -2145.1535
$ws.Range("N89").Value = -51232
$ws.Range("H134").Value = 1890.5
$ws.Range("I134").Value = 1634.1154
$ws.Range("J134").Value = 3557
$ws.Range("K134").Value = 4902.3462
$ws.Range("L134").Value = 10671
$ws.Range("M134").Value = -2367.3462
$ws.Range("N134").Value = -15741

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1456.5946
$ws.Range("I31").Value = 1292.4375
$ws.Range("J31").Value = 1581.6666
$ws.Range("K31").Value = 1292.4375
$ws.Range("L31").Value = 1581.6666
$ws.Range("M31").Value = -997.4375
$ws.Range("N31").Value = -2171.6666
$ws.Range("H34").Value = 1456.5946
$ws.Range("I34").Value = 1292.4375
$ws.Range("J34").Value = 1581.6666
$ws.Range("K34").Value = 1292.4375
$ws.Range("L34").Value = 1581.6666
$ws.Range("M34").Value = -1090.4375
$ws.Range("N34").Value = -1985.6666
$ws.Range("H105").Value = 1497.1666
$ws.Range("I105").Value = 1468.7646
$ws.Range("K105").Value = 1468.7646
$ws.Range("M105").Value = 278.2354
$ws.Range("H135").Value = 50780
$ws.Range("J135").Value = 50780
$ws.Range("L135").Value = 50780
$ws.Range("N135").Value = -60920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5829.65
$ws.Range("I3").Value = 1825
$ws.Range("J3").Value = 8499.416999999999
$ws.Range("K3").Value = 5475
$ws.Range("L3").Value = 25498.251
$ws.Range("M3").Value = -5363
$ws.Range("N3").Value = -25722.251
$ws.Range("H68").Value = 903.1111
$ws.Range("I68").Value = 471.5926
$ws.Range("J68").Value = 1226.75
$ws.Range("K68").Value = 1414.7778
$ws.Range("L68").Value = 3680.25
$ws.Range("M68").Value = -603.7778000000001
$ws.Range("N68").Value = -5302.25
$ws.Range("H71").Value = 903.1111
$ws.Range("I71").Value = 471.5926
$ws.Range("J71").Value = 1226.75
$ws.Range("K71").Value = 4244.3334
$ws.Range("L71").Value = 11040.75
$ws.Range("M71").Value = -188.3334000000004
$ws.Range("N71").Value = -19152.75
$ws.Range("H92").Value = 921.625
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 921.625
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2764.875
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -5260.875
$ws.Range("H114").Value = 1060.4348
$ws.Range("I114").Value = 629.75
$ws.Range("J114").Value = 1290.1333
$ws.Range("K114").Value = 1889.25
$ws.Range("L114").Value = 3870.3999
$ws.Range("M114").Value = 1364.75
$ws.Range("N114").Value = -10378.3999
$ws.Range("H118").Value = 6044.75
$ws.Range("I118").Value = 1679
$ws.Range("J118").Value = 7500
$ws.Range("K118").Value = 5037
$ws.Range("L118").Value = 22500
$ws.Range("M118").Value = -3794
$ws.Range("N118").Value = -24986
$ws.Range("H131").Value = 15035.267
$ws.Range("J131").Value = 1632.1528
$ws.Range("L131").Value = 4896.4584
$ws.Range("N131").Value = -14976.4584

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4181.5713
$ws.Range("I70").Value = 4100.684
$ws.Range("J70").Value = 4950
$ws.Range("K70").Value = 4100.684
$ws.Range("L70").Value = 4950
$ws.Range("M70").Value = -3830.684
$ws.Range("N70").Value = -5490
$ws.Range("H73").Value = 4181.5713
$ws.Range("I73").Value = 4100.684
$ws.Range("J73").Value = 4950
$ws.Range("K73").Value = 4100.684
$ws.Range("L73").Value = 4950
$ws.Range("M73").Value = -3164.684
$ws.Range("N73").Value = -6822
$ws.Range("H107").Value = 863.0909
$ws.Range("I107").Value = 856.2857
$ws.Range("J107").Value = 875
$ws.Range("K107").Value = 856.2857
$ws.Range("L107").Value = 875
$ws.Range("M107").Value = 1063.7143
$ws.Range("N107").Value = -4715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 783.9091
$ws.Range("I82").Value = 808.9091
$ws.Range("J82").Value = 758.9091
$ws.Range("K82").Value = 808.9091
$ws.Range("L82").Value = 758.9091
$ws.Range("M82").Value = -447.9091
$ws.Range("N82").Value = -1480.9091
$ws.Range("H85").Value = 783.9091
$ws.Range("I85").Value = 808.9091
$ws.Range("J85").Value = 758.9091
$ws.Range("K85").Value = 808.9091
$ws.Range("L85").Value = 758.9091
$ws.Range("M85").Value = 439.0909
$ws.Range("N85").Value = -3254.9091
$ws.Range("H136").Value = 3996.8948
$ws.Range("I136").Value = 1996.0667
$ws.Range("J136").Value = 11500
$ws.Range("K136").Value = 5988.2001
$ws.Range("L136").Value = 34500
$ws.Range("M136").Value = -3438.2001
$ws.Range("N136").Value = -39600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 999.6
$ws.Range("I81").Value = 999.6
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1999.2
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -938.2
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 999.6
$ws.Range("I84").Value = 999.6
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9996
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4692
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 1250
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 3750
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -8650
$ws.Range("H126").Value = 826.63635
$ws.Range("I126").Value = 823.5
$ws.Range("J126").Value = 835
$ws.Range("K126").Value = 2470.5
$ws.Range("L126").Value = 2505
$ws.Range("M126").Value = -0.5
$ws.Range("N126").Value = -7445
